# Auto-generated edit script applying numeric corrections to the
# Ixion_Profits market-data sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each sheet is an Excel Table of cached FFXIV market-board values with
# no formulas; we overwrite specific cells in H:N on the given rows with
# refreshed price/profit figures, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 10612.875
$ws.Range("I76").Value = 19801
$ws.Range("J76").Value = 5100
$ws.Range("K76").Value = 19801
$ws.Range("L76").Value = 5100
$ws.Range("M76").Value = -19486
$ws.Range("N76").Value = -5730

$ws.Range("H79").Value = 10612.875
$ws.Range("I79").Value = 19801
$ws.Range("J79").Value = 5100
$ws.Range("K79").Value = 19801
$ws.Range("L79").Value = 5100
$ws.Range("M79").Value = -18709
$ws.Range("N79").Value = -7284

$ws.Range("H86").Value = 2987.5881
$ws.Range("I86").Value = 3111.5386
$ws.Range("K86").Value = 3111.5386
$ws.Range("M86").Value = -1988.5386

$ws.Range("H89").Value = 2987.5881
$ws.Range("I89").Value = 3111.5386
$ws.Range("K89").Value = 15557.693
$ws.Range("M89").Value = -9941.692999999999

$ws.Range("H116").Value = 7208.826
$ws.Range("I116").Value = 9696.846
$ws.Range("J116").Value = 3974.4
$ws.Range("K116").Value = 9696.846
$ws.Range("L116").Value = 3974.4
$ws.Range("M116").Value = -6254.846
$ws.Range("N116").Value = -10858.4

$ws.Range("H138").Value = 1916.6102
$ws.Range("I138").Value = 1128.7941
$ws.Range("J138").Value = 2988.04
$ws.Range("K138").Value = 3386.3823
$ws.Range("L138").Value = 8964.119999999999
$ws.Range("M138").Value = 1753.6177
$ws.Range("N138").Value = -19244.12

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 8210
$ws.Range("J29").Value = 8210
$ws.Range("L29").Value = 8210
$ws.Range("N29").Value = -8826

$ws.Range("H32").Value = 1498.8687
$ws.Range("I32").Value = 1545.871
$ws.Range("K32").Value = 1545.871
$ws.Range("M32").Value = -1258.871

$ws.Range("H61").Value = 232230.25
$ws.Range("I61").Value = 6315.8096
$ws.Range("J61").Value = 438499.97
$ws.Range("K61").Value = 6315.8096
$ws.Range("L61").Value = 438499.97
$ws.Range("M61").Value = -6103.8096
$ws.Range("N61").Value = -438923.97

$ws.Range("H88").Value = 2430.3
$ws.Range("I88").Value = 2232
$ws.Range("J88").Value = 2515.2856
$ws.Range("K88").Value = 2232
$ws.Range("L88").Value = 2515.2856
$ws.Range("M88").Value = -1826
$ws.Range("N88").Value = -3327.2856

$ws.Range("H91").Value = 2430.3
$ws.Range("I91").Value = 2232
$ws.Range("J91").Value = 2515.2856
$ws.Range("K91").Value = 2232
$ws.Range("L91").Value = 2515.2856
$ws.Range("M91").Value = -828
$ws.Range("N91").Value = -5323.2856

$ws.Range("H110").Value = 2113.8823
$ws.Range("I110").Value = 1924.7142
$ws.Range("J110").Value = 2996.6667
$ws.Range("K110").Value = 1924.7142
$ws.Range("L110").Value = 2996.6667
$ws.Range("M110").Value = 120.2858000000001
$ws.Range("N110").Value = -7086.6667

$ws.Range("H136").Value = 232230.25
$ws.Range("I136").Value = 6315.8096
$ws.Range("J136").Value = 438499.97
$ws.Range("K136").Value = 18947.4288
$ws.Range("L136").Value = 1315499.91
$ws.Range("M136").Value = -16397.4288
$ws.Range("N136").Value = -1320599.91

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 4500
$ws.Range("I16").Value = 4500
$ws.Range("K16").Value = 4500
$ws.Range("M16").Value = -4330

$ws.Range("H20").Value = 19449.143
$ws.Range("J20").Value = 125200
$ws.Range("L20").Value = 125200
$ws.Range("N20").Value = -125694

$ws.Range("H94").Value = 1434.7916
$ws.Range("I94").Value = 972.3333
$ws.Range("J94").Value = 2205.5557
$ws.Range("K94").Value = 972.3333
$ws.Range("L94").Value = 2205.5557
$ws.Range("M94").Value = -521.3333
$ws.Range("N94").Value = -3107.5557

$ws.Range("H134").Value = 36755.188
$ws.Range("I134").Value = 7167.6
$ws.Range("J134").Value = 86067.836
$ws.Range("K134").Value = 21502.8
$ws.Range("L134").Value = 258203.508
$ws.Range("M134").Value = -18967.8
$ws.Range("N134").Value = -263273.508

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 500006000

$ws.Range("H27").Value = 500006000

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 792.625
$ws.Range("I92").Value = 748.6667
$ws.Range("J92").Value = 800.3823
$ws.Range("K92").Value = 2246.0001
$ws.Range("L92").Value = 2401.1469
$ws.Range("M92").Value = -998.0001000000002
$ws.Range("N92").Value = -4897.1469

$ws.Range("H107").Value = 428.2069
$ws.Range("I107").Value = 438.58334
$ws.Range("K107").Value = 1315.75002
$ws.Range("M107").Value = 604.2499800000001

$ws.Range("H122").Value = 4585
$ws.Range("I122").Value = 410.1111
$ws.Range("J122").Value = 12934.777
$ws.Range("K122").Value = 3690.9999
$ws.Range("L122").Value = 116412.993
$ws.Range("M122").Value = -1240.9999
$ws.Range("N122").Value = -121312.993

$ws.Range("H131").Value = 2084224.8
$ws.Range("I131").Value = 5263688
$ws.Range("J131").Value = 1128.2759
$ws.Range("K131").Value = 15791064
$ws.Range("L131").Value = 3384.8277
$ws.Range("M131").Value = -15786024
$ws.Range("N131").Value = -13464.8277

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6200.4287
$ws.Range("I70").Value = 6700.6665
$ws.Range("J70").Value = 5300
$ws.Range("K70").Value = 6700.6665
$ws.Range("L70").Value = 5300
$ws.Range("M70").Value = -6430.6665
$ws.Range("N70").Value = -5840

$ws.Range("H73").Value = 6200.4287
$ws.Range("I73").Value = 6700.6665
$ws.Range("J73").Value = 5300
$ws.Range("K73").Value = 6700.6665
$ws.Range("L73").Value = 5300
$ws.Range("M73").Value = -5764.6665
$ws.Range("N73").Value = -7172

$ws.Range("H136").Value = 4969.143
$ws.Range("J136").Value = 4969.143
$ws.Range("L136").Value = 14907.429
$ws.Range("N136").Value = -20007.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 1800
$ws.Range("J21").Value = 2000
$ws.Range("L21").Value = 2000
$ws.Range("N21").Value = -2348

$ws.Range("H93").Value = 1059.8
$ws.Range("I93").Value = 1014
$ws.Range("J93").Value = 1166.6666
$ws.Range("K93").Value = 1014
$ws.Range("L93").Value = 1166.6666
$ws.Range("M93").Value = 234
$ws.Range("N93").Value = -3662.6666

$ws.Range("H122").Value = 1735483.9
$ws.Range("I122").Value = 2307648.8
$ws.Range("J122").Value = 626914.3
$ws.Range("K122").Value = 6922946.399999999
$ws.Range("L122").Value = 1880742.9
$ws.Range("M122").Value = -6920496.399999999
$ws.Range("N122").Value = -1885642.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 20000
$ws.Range("J64").Value = 20000
$ws.Range("L64").Value = 20000
$ws.Range("N64").Value = -20496

$ws.Range("H67").Value = 20000
$ws.Range("J67").Value = 20000
$ws.Range("L67").Value = 20000
$ws.Range("N67").Value = -21716

